# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Comercializadora del Agro de Limarí - Limón"
# right above the existing row 308, pushing all subsequent rows down by 3.
# (Net result: dimension grows from A1:T424 to A1:T427.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 308..424 down to 311..427 by inserting 3 blank rows at 308.
$ws.Rows("308:310").Insert()

# Common (constant-across-sheet) field values shared by every data row.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad  = "Sin especificar"
$unidad    = "$/malla 16 kilos"
$origen    = "Provincia de Limarí"

# New row data: Fecha(D), Calidad(L), Volumen(M), Precio min(N), Precio max(O),
# Precio prom pond(P), Precio $/Kg(S), Kg/unidad(T).
$newRows = @(
    @{ Row = 308; Fecha = 44510; Calidad = "1a amarillo"; Volumen = 900; PMin = 5300; PMax = 5500; PProm = 5400; PKg = 338;  Kg = 16 },
    @{ Row = 309; Fecha = 44510; Calidad = "2a amarillo"; Volumen = 750; PMin = 3800; PMax = 4000; PProm = 3900; PKg = 244;  Kg = 16 },
    @{ Row = 310; Fecha = 44510; Calidad = "3a amarillo"; Volumen = 470; PMin = 2800; PMax = 3000; PProm = 2900; PKg = 181;  Kg = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.Kg
}
